# Update average_county_temperature (column AA) values for specific facility rows
# to reflect updated NOAA temperature data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 15.74228395061728
    3  = 15.74228395061728
    4  = 15.74228395061728
    5  = 15.74228395061728
    6  = 15.74228395061728
    7  = 15.74228395061728
    8  = 13.46442495126706
    9  = 13.46442495126706
    10 = 18.89814814814816
    11 = 18.89814814814816
    12 = 13.46442495126706
    13 = 13.46442495126706
    20 = 15.74228395061728
    21 = 15.74228395061728
    22 = 21.28240740740739
    23 = 21.28240740740739
    24 = 15.74228395061728
    25 = 15.74228395061728
}

foreach ($row in $updates.Keys) {
    $ws.Range("AA$row").Value = $updates[$row]
}
